$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out columns B:D (Username/password/placeholder table) ---
$ws.Columns("B:D").Delete()

# --- Re-key column A as a simple numbered "command history" list ---
# A1 keeps the existing "Facebook" header (shared string / style already s="1").
# A2 already holds 1 / style 1 -- leave as-is, just normalize row height below.
# A3 already holds 3 / style 1 -- leave as-is too.

$values = @(1,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,24)

for ($i = 0; $i -lt $values.Length; $i++) {
    $r = $i + 2
    $target = $ws.Cells.Item($r, 1)
    $target.Value = $values[$i]
    # Copy column-A's existing number style (s="1") onto every row so newly
    # created rows 4-38 match the look of the pre-existing rows 2-3.
    $ws.Range("A2").Copy($target)
    $target.Value = $values[$i]
}

# Normalize row heights back to the sheet default (drops the stray ht="13"
# left over on rows 2/3 from the old table layout).
$ws.Rows("1:38").AutoFit()

# Column A width to match the new narrow single-column layout.
$ws.Columns.Item(1).ColumnWidth = 13.6666666666667

# Selection moves to D19 (matches the new "history" selection/cursor spot).
$ws.Range("D19").Select()
